# Updated coefficients from waste experts
$wb = $excel.ActiveWorkbook

# --- Sheet "default" -------------------------------------------------
$ws1 = $wb.Worksheets.Item("default")
$ws1.Range("C3").Value = 0.4
$ws1.Range("I3").Formula = "=0.014/3"
$ws1.Range("I4").Value = 0.0045
$ws1.Range("E5").Value = 0.54
$ws1.Range("I5").Formula = "=0.014/3"

# --- Sheet "EP_HH" -----------------------------------------------------
$ws2 = $wb.Worksheets.Item("EP_HH")
$ws2.Range("I10").Value = 0.0045
$ws2.Range("I13").Formula = "=0.014/3"
$ws2.Activate()
$ws2.Range("I13").Select()

# --- Sheet "A" -----------------------------------------------------
$ws3 = $wb.Worksheets.Item("A")
$ws3.Range("C3").Value = 0.4
$ws3.Range("I3").Formula = "=0.014/3"
$ws3.Range("I4").Value = 0.0045
$ws3.Range("E5").Value = 0.54
$ws3.Range("I5").Formula = "=0.014/3"
$ws3.Activate()
$ws3.Range("I4").Select()

# --- Sheet "C10-C12" -------------------------------------------------
$ws4 = $wb.Worksheets.Item("C10-C12")
$ws4.Range("C3").Value = 0.4
$ws4.Range("I3").Formula = "=0.014/3"
$ws4.Range("I4").Value = 0.0045
$ws4.Range("E5").Value = 0.54
$ws4.Range("I5").Formula = "=0.014/3"
$ws4.Activate()
$ws4.Range("E5").Select()

# --- Sheet "G-U_X_G4677" ----------------------------------------------
$ws5 = $wb.Worksheets.Item("G-U_X_G4677")
$ws5.Range("C3").Value = 0.4
$ws5.Range("I3").Formula = "=0.014/3"
$ws5.Range("I4").Value = 0.0045
$ws5.Range("E5").Value = 0.54
$ws5.Range("I5").Formula = "=0.014/3"
$ws5.Activate()
$ws5.Range("E6").Select()

# Final active sheet/selection: "default" sheet, cell I5
$ws1.Activate()
$ws1.Range("I5").Select()
